# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures for every coin row (2-51) to match the latest scrape, as
# produced by the scheduled GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as '39.589.22' or '227.21'. Some of
# those look like plain numbers to Excel's automatic type detection, so
# the whole column is temporarily forced to Text format while the new
# values are written, then restored to the default (Normal) style so the
# cells keep the original 'no explicit style' look of column D.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '39.589.22'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '2.168.71'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '227.21'
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("D6").Value = '0.622'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").Value = '62.87'
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.390'
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("D10").Value = '0.0846'
$ws.Range("E10").Value = '  -0.65%  '
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").Value = '15.87'
$ws.Range("E12").Value = '  -0.93%  '
$ws.Range("D13").Value = '2.489.81'
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("D14").Value = '21.71'
$ws.Range("E14").Value = '  -2.36%  '
$ws.Range("D15").Value = '0.806'
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("D16").Value = '5.46'
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D17").Value = '2.162.23'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").Value = '39.606.79'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").Value = '71.64'
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("D20").Value = '0.0₃0889'
$ws.Range("E20").Value = '  +4.42%  '
$ws.Range("D21").Value = '6.00'
$ws.Range("E21").Value = '  -2.29%  '
$ws.Range("D22").Value = '227.81'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '2.35'
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").Value = '2.32'
$ws.Range("E25").Value = '  -4.04%  '
$ws.Range("D26").Value = '170.37'
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").Value = '9.45'
$ws.Range("E27").Value = '  -3.06%  '
$ws.Range("D28").Value = '0.137'
$ws.Range("E28").Value = '  -0.82%  '
$ws.Range("D29").Value = '1.44'
$ws.Range("E29").Value = '  +2.45%  '
$ws.Range("D30").Value = '19.67'
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("D31").Value = '2.67'
$ws.Range("E31").Value = '  +4.23%  '
$ws.Range("D32").Value = '0.122'
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("D33").Value = '4.48'
$ws.Range("E33").Value = '  -2.85%  '
$ws.Range("D34").Value = '4.70'
$ws.Range("E34").Value = '  -2.52%  '
$ws.Range("D35").Value = '6.97'
$ws.Range("E35").Value = '  -2.29%  '
$ws.Range("D36").Value = '0.0616'
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("D37").Value = '3.78'
$ws.Range("E37").Value = '  +6.70%  '
$ws.Range("D38").Value = '2.39'
$ws.Range("E38").Value = '  -0.62%  '
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").Value = '4.91'
$ws.Range("E40").Value = '  +17.75%  '
$ws.Range("D41").Value = '102.26'
$ws.Range("E41").Value = '  -0.58%  '
$ws.Range("D42").Value = '0.0227'
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").Value = '17.74'
$ws.Range("E43").Value = '  -2.30%  '
$ws.Range("D44").Value = '1.513.31'
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("D45").Value = '1.20'
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("D46").Value = '7.89'
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("D47").Value = '2.80'
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("D48").Value = '0.0914'
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("D49").Value = '1.09'
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("D50").Value = '0.000196'
$ws.Range("E50").Value = '  +31.39%  '
$ws.Range("D51").Value = '2.371.82'
$ws.Range("E51").Value = '  +0.74%  '

# Restore the default cell style now that the text values are locked in.
$dRange.Style = "Normal"
